# Rename document/table attributes to lowerCamelCase in the ObjTables
# header rows embedded as inline strings in cell A1/A2 of each sheet.

$wb = $excel.ActiveWorkbook

$tocSheet = $wb.Worksheets.Item("!!_Table of contents")
$tocSheet.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$tocSheet.Range("A2").Value = "!!ObjTables type='TableOfContents' description='Table of contents' date='2019-09-24 17:36:04' objTablesVersion='0.0.8'"

$schemaSheet = $wb.Worksheets.Item("!!_Schema")
$schemaSheet.Range("A1").Value = "!!ObjTables type='Schema' description='Table/model and column/attribute definitions' date='2019-09-24 17:36:04' objTablesVersion='0.0.8'"

$transactionSheet = $wb.Worksheets.Item("!!Transaction")
$transactionSheet.Range("A1").Value = "!!ObjTables type='Data' id='Transaction' description='Stores transactions' name='Transaction' date='2019-09-24 17:36:04' objTablesVersion='0.0.8'"
